$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Players")
$ws2 = $wb.Worksheets.Item("OwnerTotals")

# Column G width 18 -> 17 on Players sheet (column 7)
# NOTE: Excel's ColumnWidth (character units) round-trips to the stored
# OOXML "width" with a fixed +5/6 offset for this workbook's default font,
# so we back-solve the character width that serializes to exactly 17.
$ws1.Columns.Item(7).ColumnWidth = 16.166666666666668

# Row 2
$ws1.Range("G2").Value = "Halftime"

# Row 3
$ws1.Range("G3").Value = "3:47 - 2nd Half"
$ws1.Range("H3").Value = 12
$ws1.Range("J3").Value = 6
$ws1.Range("O3").Value = 32

# Row 4
$ws1.Range("G4").Value = "3:47 - 2nd Half"
$ws1.Range("H4").Value = 10
$ws1.Range("J4").Value = 3
$ws1.Range("O4").Value = 17

# Row 5
$ws1.Range("G5").Value = "Halftime"

# Row 6
$ws1.Range("G6").Value = "3:47 - 2nd Half"
$ws1.Range("O6").Value = 34

# Row 7
$ws1.Range("G7").Value = "3:47 - 2nd Half"
$ws1.Range("H7").Value = 19
$ws1.Range("I7").Value = 22
$ws1.Range("J7").Value = 2
$ws1.Range("O7").Value = 31

# Row 8
$ws1.Range("G8").Value = "3:47 - 2nd Half"
$ws1.Range("O8").Value = 31

# Row 9
$ws1.Range("G9").Value = "Halftime"

# Row 10
$ws1.Range("G10").Value = "Halftime"

# Row 11
$ws1.Range("G11").Value = "3:47 - 2nd Half"

# Row 12
$ws1.Range("G12").Value = "Halftime"

# Row 13
$ws1.Range("G13").Value = "Halftime"

# Row 14
$ws1.Range("G14").Value = "3:47 - 2nd Half"

# Row 15
$ws1.Range("G15").Value = "3:47 - 2nd Half"

# Row 16
$ws1.Range("G16").Value = "Halftime"

# Row 17
$ws1.Range("G17").Value = "Halftime"

# Row 18
$ws1.Range("G18").Value = "3:47 - 2nd Half"
$ws1.Range("H18").Value = 19
$ws1.Range("O18").Value = 35

# Row 19
$ws1.Range("G19").Value = "3:47 - 2nd Half"
$ws1.Range("H19").Value = 14
$ws1.Range("I19").Value = 13
$ws1.Range("K19").Value = 3
$ws1.Range("O19").Value = 37

# Row 20
$ws1.Range("G20").Value = "3:47 - 2nd Half"

# Row 21
$ws1.Range("G21").Value = "3:47 - 2nd Half"
$ws1.Range("H21").Value = 17
$ws1.Range("I21").Value = 9
$ws1.Range("J21").Value = 11
$ws1.Range("O21").Value = 27

# Row 22
$ws1.Range("G22").Value = "Halftime"

# Row 23
$ws1.Range("G23").Value = "3:47 - 2nd Half"
$ws1.Range("H23").Value = 10
$ws1.Range("O23").Value = 34

# Row 24
$ws1.Range("G24").Value = "Halftime"

# Row 25
$ws1.Range("G25").Value = "3:47 - 2nd Half"

# Row 26
$ws1.Range("G26").Value = "3:47 - 2nd Half"
$ws1.Range("H26").Value = 15
$ws1.Range("I26").Value = 12
$ws1.Range("J26").Value = 4
$ws1.Range("O26").Value = 18

# Row 27
$ws1.Range("G27").Value = "Halftime"

# Row 28
$ws1.Range("G28").Value = "Halftime"

# Row 29
$ws1.Range("G29").Value = "Halftime"

# Row 30
$ws1.Range("G30").Value = "Halftime"

# Row 31
$ws1.Range("G31").Value = "Halftime"

# Row 32
$ws1.Range("D32").Value = "Nicholas Randall"
$ws1.Range("E32").Value = "MIZ"
$ws1.Range("F32").Value = "MIZ@ALA"
$ws1.Range("G32").Value = "Halftime"
$ws1.Range("J32").Value = 0
$ws1.Range("N32").Value = 0
$ws1.Range("O32").Value = 1

# Row 33
$ws1.Range("D33").Value = "D.J. Wagner"
$ws1.Range("E33").Value = "ARK"
$ws1.Range("F33").Value = "ARK@OU"
$ws1.Range("G33").Value = "3:47 - 2nd Half"
$ws1.Range("H33").Value = -1
$ws1.Range("J33").Value = 2
$ws1.Range("N33").Value = 1
$ws1.Range("O33").Value = 13

# Row 34
$ws1.Range("G34").Value = "Halftime"

# OwnerTotals sheet updates
$ws2.Range("B2").Value = 19
$ws2.Range("B3").Value = 17
